$d = $word.ActiveDocument
$r = $d.Content
$found1 = $r.Find.Execute("<digitalSignature", $false, $false, $false, $false, $false, $true, 1, $false, "<digitalSignature", 2)
Write-Output ("Found1: " + $found1)

$r2 = $d.Content
$found2 = $r2.Find.Execute("Buyer>", $false, $false, $false, $false, $false, $true, 1, $false, "Buyer>", 2)
Write-Output ("Found2: " + $found2)

# positions: full para range 59-84, text "<digitalSignature_Buyer>" len 25, ">" at the very end = position 83 (84 is end)
$sub = $d.Range(83, 84)
Write-Output ("sub text: [" + $sub.Text + "]")
$sub.Font.Bold = $true
$sub.Font.Bold = $false
